$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: convert A38 and D38 from text to numbers (B38/C38 stay as-is)
$ws.Cells.Item(38, 1).Value = 31
$ws.Cells.Item(38, 4).Value = 11

# Rows 39-90: new data rows (A,D numeric; B,C text)
$newRows = @(
    @(39, 29, "Light Rain", "08/12/2024", 13),
    @(40, 29, "Rain Shower", "08/12/2024", 14),
    @(41, 29, "Rain Shower", "08/12/2024", 14),
    @(42, 31, "Mostly Cloudy", "08/12/2024", 14),
    @(43, 31, "Mostly Cloudy", "08/12/2024", 14),
    @(44, 29, "Light Rain", "08/12/2024", 14),
    @(45, 27, "Rain", "08/12/2024", 14),
    @(46, 29, "Light Rain", "08/12/2024", 15),
    @(47, 29, "Light Rain", "08/12/2024", 15),
    @(48, 29, "Light Rain", "08/12/2024", 15),
    @(49, 29, "Light Rain", "08/12/2024", 15),
    @(50, 29, "Light Rain", "08/12/2024", 15),
    @(51, 30, "Mostly Cloudy", "08/12/2024", 16),
    @(52, 27, "Rain", "08/12/2024", 16),
    @(53, 28, "Light Rain", "08/12/2024", 16),
    @(54, 28, "Light Rain", "08/12/2024", 16),
    @(55, 28, "Rain Shower", "08/12/2024", 16),
    @(56, 30, "Mostly Cloudy", "08/12/2024", 16),
    @(57, 29, "Mostly Cloudy", "08/12/2024", 17),
    @(58, 29, "Mostly Cloudy", "08/12/2024", 17),
    @(59, 28, "Rain Shower", "08/12/2024", 17),
    @(60, 28, "Rain Shower", "08/12/2024", 17),
    @(61, 27, "Rain", "08/12/2024", 17),
    @(62, 29, "Mostly Cloudy", "08/12/2024", 17),
    @(63, 29, "Mostly Cloudy", "08/12/2024", 18),
    @(64, 27, "Partly Cloudy", "08/12/2024", 21),
    @(65, 27, "Partly Cloudy", "08/12/2024", 22),
    @(66, 27, "Partly Cloudy", "08/12/2024", 22),
    @(67, 27, "Partly Cloudy", "08/12/2024", 22),
    @(68, 27, "Partly Cloudy", "08/12/2024", 22),
    @(69, 27, "Partly Cloudy", "08/12/2024", 22),
    @(70, 27, "Partly Cloudy", "08/12/2024", 23),
    @(71, 27, "Partly Cloudy", "08/12/2024", 23),
    @(72, 27, "Partly Cloudy", "08/12/2024", 23),
    @(73, 27, "Partly Cloudy", "08/12/2024", 23),
    @(74, 27, "Partly Cloudy", "08/12/2024", 23),
    @(75, 27, "Partly Cloudy", "08/12/2024", 23),
    @(76, 27, "Partly Cloudy", "08/13/2024", 0),
    @(77, 27, "Partly Cloudy", "08/13/2024", 0),
    @(78, 27, "Partly Cloudy", "08/13/2024", 0),
    @(79, 27, "Partly Cloudy", "08/13/2024", 0),
    @(80, 27, "Partly Cloudy", "08/13/2024", 0),
    @(81, 27, "Partly Cloudy", "08/13/2024", 1),
    @(82, 29, "Cloudy", "08/13/2024", 8),
    @(83, 29, "Cloudy", "08/13/2024", 8),
    @(84, 29, "Cloudy", "08/13/2024", 8),
    @(85, 29, "Cloudy", "08/13/2024", 8),
    @(86, 30, "Cloudy", "08/13/2024", 8),
    @(87, 30, "Mostly Cloudy", "08/13/2024", 9),
    @(88, 30, "Mostly Cloudy", "08/13/2024", 9),
    @(89, 31, "Mostly Cloudy", "08/13/2024", 9),
    @(90, 31, "Mostly Cloudy", "08/13/2024", 9)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = "'" + $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Row 91: new row, all text (A91="31" text, D91="09" text with leading zero)
$ws.Cells.Item(91, 1).Value = "'31"
$ws.Cells.Item(91, 2).Value = "Mostly Cloudy"
$ws.Cells.Item(91, 3).Value = "'08/13/2024"
$ws.Cells.Item(91, 4).Value = "'09"

Write-Output "done"
